$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number-looking string (e.g. "216.95").
# Excel would auto-convert such text to a numeric value on assignment, but the
# source data stores these as literal text (inlineStr) - force text formatting
# first so the written value keeps its exact string representation (incl. any
# trailing zero, like "146.70") instead of becoming a float.
$textForceCells = @('D5', 'D6', 'D9', 'D10', 'D11', 'D13', 'D16', 'D19', 'D21', 'D25', 'D39', 'D44', 'D46', 'D50', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data scraped by the GitHub Actions job.
$ws.Range('D2').Value = '26.874.49'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.638.02'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.66%  '
$ws.Range('D5').Value = '216.95'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = '0.511'
$ws.Range('E6').Value = '  +1.87%  '
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').Value = '  +1.70%  '
$ws.Range('D9').Value = '0.0625'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').Value = '19.91'
$ws.Range('E10').Value = '  +3.84%  '
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').Value = '1.866.98'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.11'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.601.26'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').Value = '66.99'
$ws.Range('E16').Value = '  +3.00%  '
$ws.Range('D17').Value = '26.871.81'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').Value = '219.42'
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '6.84'
$ws.Range('E21').Value = '  +3.94%  '
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '146.70'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').Value = '1.257.86'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = '0.834'
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('D43').Value = '1.777.43'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '61.81'
$ws.Range('E44').Value = '  +1.52%  '
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('D46').Value = '91.55'
$ws.Range('E46').Value = '  -0.87%  '
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('E48').Value = '  +3.27%  '
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('D50').Value = '7.65'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('D51').Value = '0.0961'
$ws.Range('E51').Value = '  -0.15%  '
